# Insert two new data rows at rows 95-96 (this shifts all existing rows
# 95..148 down to 97..150, and Excel auto-updates the sheet dimension).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("95:96").Insert()

# Row 95: new "Especial" record for Mango from Brasil
$ws.Cells.Item(95, 1).Value = 1
$ws.Cells.Item(95, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(95, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(95, 4).Value = 44818
$ws.Cells.Item(95, 5).Value = 15
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100108
$ws.Cells.Item(95, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(95, 9).Value = 100108002
$ws.Cells.Item(95, 10).Value = "Mango"
$ws.Cells.Item(95, 11).Value = "Sin especificar"
$ws.Cells.Item(95, 12).Value = "Especial"
$ws.Cells.Item(95, 13).Value = 570
$ws.Cells.Item(95, 14).Value = 7500
$ws.Cells.Item(95, 15).Value = 8000
$ws.Cells.Item(95, 16).Value = 7750
$ws.Cells.Item(95, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(95, 18).Value = "Brasil"
$ws.Cells.Item(95, 19).Value = 1938
$ws.Cells.Item(95, 20).Value = 4

# Row 96: new "Primera" record for Mango from Brasil
$ws.Cells.Item(96, 1).Value = 1
$ws.Cells.Item(96, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(96, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(96, 4).Value = 44818
$ws.Cells.Item(96, 5).Value = 15
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100108
$ws.Cells.Item(96, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(96, 9).Value = 100108002
$ws.Cells.Item(96, 10).Value = "Mango"
$ws.Cells.Item(96, 11).Value = "Sin especificar"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 570
$ws.Cells.Item(96, 14).Value = 7500
$ws.Cells.Item(96, 15).Value = 8000
$ws.Cells.Item(96, 16).Value = 7750
$ws.Cells.Item(96, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(96, 18).Value = "Brasil"
$ws.Cells.Item(96, 19).Value = 1938
$ws.Cells.Item(96, 20).Value = 4
